$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1084441.6
$ws.Range("I17").Value = 398.92307
$ws.Range("J17").Value = 1771883.4
$ws.Range("K17").Value = 1196.76921
$ws.Range("L17").Value = 5315650.199999999
$ws.Range("M17").Value = -1028.76921
$ws.Range("N17").Value = -5315986.199999999
$ws.Range("H132").Value = 3194.925
$ws.Range("I132").Value = 3059.125
$ws.Range("J132").Value = 3738.125
$ws.Range("K132").Value = 9177.375
$ws.Range("L132").Value = 11214.375
$ws.Range("M132").Value = -6647.375
$ws.Range("N132").Value = -16274.375
$ws.Range("H135").Value = 808.88
$ws.Range("I135").Value = 459.45
$ws.Range("J135").Value = 2206.6
$ws.Range("K135").Value = 4135.05
$ws.Range("L135").Value = 19859.4
$ws.Range("M135").Value = -1600.05
$ws.Range("N135").Value = -24929.4
$ws.Range("H137").Value = 1154.2667
$ws.Range("I137").Value = 990.55316
$ws.Range("J137").Value = 1429.0714
$ws.Range("K137").Value = 2971.65948
$ws.Range("L137").Value = 4287.2142
$ws.Range("M137").Value = -421.6594800000003
$ws.Range("N137").Value = -9387.2142
$ws.Range("H138").Value = 1395.66
$ws.Range("I138").Value = 642.7838
$ws.Range("J138").Value = 1837.8254
$ws.Range("K138").Value = 1928.3514
$ws.Range("L138").Value = 5513.4762
$ws.Range("M138").Value = 3211.6486
$ws.Range("N138").Value = -15793.4762

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 671637.0600000001
$ws.Range("I32").Value = 767487.4399999999
$ws.Range("K32").Value = 767487.4399999999
$ws.Range("M32").Value = -767200.4399999999
$ws.Range("H61").Value = 7754139.5
$ws.Range("I61").Value = 9805630
$ws.Range("J61").Value = 4066.3333
$ws.Range("K61").Value = 9805630
$ws.Range("L61").Value = 4066.3333
$ws.Range("M61").Value = -9805418
$ws.Range("N61").Value = -4490.3333
$ws.Range("H74").Value = 1589.46
$ws.Range("I74").Value = 865.08
$ws.Range("J74").Value = 2313.84
$ws.Range("K74").Value = 865.08
$ws.Range("L74").Value = 2313.84
$ws.Range("M74").Value = 8.919999999999959
$ws.Range("N74").Value = -4061.84
$ws.Range("H77").Value = 1589.46
$ws.Range("I77").Value = 865.08
$ws.Range("J77").Value = 2313.84
$ws.Range("K77").Value = 4325.400000000001
$ws.Range("L77").Value = 11569.2
$ws.Range("M77").Value = 42.59999999999945
$ws.Range("N77").Value = -20305.2
$ws.Range("H132").Value = 2803.9302
$ws.Range("I132").Value = 2254.2258
$ws.Range("J132").Value = 4224
$ws.Range("K132").Value = 6762.6774
$ws.Range("L132").Value = 12672
$ws.Range("M132").Value = -4232.6774
$ws.Range("N132").Value = -17732
$ws.Range("H136").Value = 7754139.5
$ws.Range("I136").Value = 9805630
$ws.Range("J136").Value = 4066.3333
$ws.Range("K136").Value = 29416890
$ws.Range("L136").Value = 12198.9999
$ws.Range("M136").Value = -29414340
$ws.Range("N136").Value = -17298.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4316.7583
$ws.Range("I31").Value = 1343.3871
$ws.Range("K31").Value = 1343.3871
$ws.Range("M31").Value = -1048.3871
$ws.Range("H34").Value = 4316.7583
$ws.Range("I34").Value = 1343.3871
$ws.Range("K34").Value = 1343.3871
$ws.Range("M34").Value = -1141.3871
$ws.Range("H58").Value = 766.8357999999999
$ws.Range("I58").Value = 519.2083
$ws.Range("J58").Value = 1392.421
$ws.Range("K58").Value = 519.2083
$ws.Range("L58").Value = 1392.421
$ws.Range("M58").Value = -316.2083
$ws.Range("N58").Value = -1798.421
$ws.Range("H132").Value = 4168126.5
$ws.Range("I132").Value = 1142.5927
$ws.Range("J132").Value = 12822631
$ws.Range("K132").Value = 3427.7781
$ws.Range("L132").Value = 38467893
$ws.Range("M132").Value = -897.7780999999995
$ws.Range("N132").Value = -38472953
$ws.Range("H134").Value = 4019.6511
$ws.Range("I134").Value = 4249.2812
$ws.Range("J134").Value = 3351.6365
$ws.Range("K134").Value = 12747.8436
$ws.Range("L134").Value = 10054.9095
$ws.Range("M134").Value = -10212.8436
$ws.Range("N134").Value = -15124.9095
$ws.Range("H136").Value = 766.8357999999999
$ws.Range("I136").Value = 519.2083
$ws.Range("J136").Value = 1392.421
$ws.Range("K136").Value = 1557.6249
$ws.Range("L136").Value = 4177.263
$ws.Range("M136").Value = 992.3751
$ws.Range("N136").Value = -9277.262999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1040.159
$ws.Range("I5").Value = 272.85184
$ws.Range("K5").Value = 818.5555199999999
$ws.Range("M5").Value = -706.5555199999999
$ws.Range("H92").Value = 900
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H131").Value = 2881.4915
$ws.Range("I131").Value = 361.6111
$ws.Range("J131").Value = 3987.7805
$ws.Range("K131").Value = 1084.8333
$ws.Range("L131").Value = 11963.3415
$ws.Range("M131").Value = 3955.1667
$ws.Range("N131").Value = -22043.3415
$ws.Range("H133").Value = 13800.154
$ws.Range("I133").Value = 8321.666999999999
$ws.Range("J133").Value = 18496
$ws.Range("K133").Value = 24965.001
$ws.Range("L133").Value = 55488
$ws.Range("M133").Value = -19905.001
$ws.Range("N133").Value = -65608
$ws.Range("H135").Value = 1040.159
$ws.Range("I135").Value = 272.85184
$ws.Range("K135").Value = 2455.66656
$ws.Range("M135").Value = 79.33344000000034

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3203.3684
$ws.Range("I132").Value = 2661.5925
$ws.Range("J132").Value = 4533.1816
$ws.Range("K132").Value = 7984.7775
$ws.Range("L132").Value = 13599.5448
$ws.Range("M132").Value = -5454.7775
$ws.Range("N132").Value = -18659.5448

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 62503710
$ws.Range("I82").Value = 100003944
$ws.Range("J82").Value = 3335.3333
$ws.Range("K82").Value = 100003944
$ws.Range("L82").Value = 3335.3333
$ws.Range("M82").Value = -100003583
$ws.Range("N82").Value = -4057.3333
$ws.Range("H85").Value = 62503710
$ws.Range("I85").Value = 100003944
$ws.Range("J85").Value = 3335.3333
$ws.Range("K85").Value = 100003944
$ws.Range("L85").Value = 3335.3333
$ws.Range("M85").Value = -100002696
$ws.Range("N85").Value = -5831.3333
$ws.Range("H132").Value = 2363.3433
$ws.Range("I132").Value = 2015.0358
$ws.Range("J132").Value = 4136.5454
$ws.Range("K132").Value = 6045.107400000001
$ws.Range("L132").Value = 12409.6362
$ws.Range("M132").Value = -3515.107400000001
$ws.Range("N132").Value = -17469.6362
$ws.Range("H136").Value = 2733482.5
$ws.Range("I136").Value = 1011.1818
$ws.Range("J136").Value = 9805761
$ws.Range("K136").Value = 3033.5454
$ws.Range("L136").Value = 29417283
$ws.Range("M136").Value = -483.5454
$ws.Range("N136").Value = -29422383

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2603.75
$ws.Range("I81").Value = 2972.5
$ws.Range("J81").Value = 1866.25
$ws.Range("K81").Value = 5945
$ws.Range("L81").Value = 3732.5
$ws.Range("M81").Value = -4884
$ws.Range("N81").Value = -5854.5
$ws.Range("H84").Value = 2603.75
$ws.Range("I84").Value = 2972.5
$ws.Range("J84").Value = 1866.25
$ws.Range("K84").Value = 29725
$ws.Range("L84").Value = 18662.5
$ws.Range("M84").Value = -24421
$ws.Range("N84").Value = -29270.5
$ws.Range("H122").Value = 2707.0217
$ws.Range("I122").Value = 2501.4546
$ws.Range("J122").Value = 3228.8462
$ws.Range("K122").Value = 7504.3638
$ws.Range("L122").Value = 9686.5386
$ws.Range("M122").Value = -5054.3638
$ws.Range("N122").Value = -14586.5386
$ws.Range("H132").Value = 4558793
$ws.Range("I132").Value = 1473.7954
$ws.Range("J132").Value = 14584895
$ws.Range("K132").Value = 4421.3862
$ws.Range("L132").Value = 43754685
$ws.Range("M132").Value = -1891.3862
$ws.Range("N132").Value = -43759745
$ws.Range("H136").Value = 1450.54
$ws.Range("I136").Value = 1412.421
$ws.Range("J136").Value = 1571.25
$ws.Range("K136").Value = 4237.263
$ws.Range("L136").Value = 4713.75
$ws.Range("M136").Value = -1687.263
$ws.Range("N136").Value = -9813.75

# --- Remove N92 cell entirely (CUL) to match target OOXML (no leftover empty cell) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N92").ClearContents()
